$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 739
$ws.Range("I33").Value = 93.14286
$ws.Range("K33").Value = 93.14286
$ws.Range("M33").Value = 135.85714
$ws.Range("H87").Value = 44000
$ws.Range("J87").Value = 44000
$ws.Range("L87").Value = 44000
$ws.Range("N87").Value = -46496
$ws.Range("H88").Value = 24509
$ws.Range("I88").Value = 6598.3335
$ws.Range("J88").Value = 33464.332
$ws.Range("K88").Value = 6598.3335
$ws.Range("L88").Value = 33464.332
$ws.Range("M88").Value = -6192.3335
$ws.Range("N88").Value = -34276.332
$ws.Range("H90").Value = 44000
$ws.Range("J90").Value = 44000
$ws.Range("L90").Value = 132000
$ws.Range("N90").Value = -144480
$ws.Range("H91").Value = 24509
$ws.Range("I91").Value = 6598.3335
$ws.Range("J91").Value = 33464.332
$ws.Range("K91").Value = 6598.3335
$ws.Range("L91").Value = 33464.332
$ws.Range("M91").Value = -5194.3335
$ws.Range("N91").Value = -36272.332
$ws.Range("H115").Value = 373.5
$ws.Range("I115").Value = 373.5
$ws.Range("K115").Value = 1120.5
$ws.Range("M115").Value = 446.5
$ws.Range("H125").Value = 103160.7
$ws.Range("J125").Value = 337941.34
$ws.Range("L125").Value = 3041472.06
$ws.Range("N125").Value = -3046392.06
$ws.Range("H137").Value = 1696.55
$ws.Range("I137").Value = 1575.0526
$ws.Range("K137").Value = 4725.1578
$ws.Range("M137").Value = -2175.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 3794.5
$ws.Range("J46").Value = 3794.5
$ws.Range("L46").Value = 3794.5
$ws.Range("N46").Value = -4432.5
$ws.Range("H55").Value = 28333.334
$ws.Range("J55").Value = 28333.334
$ws.Range("L55").Value = 28333.334
$ws.Range("N55").Value = -28963.334
$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -41996
$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -129984
$ws.Range("H110").Value = 504.66666
$ws.Range("I110").Value = 389.5
$ws.Range("K110").Value = 389.5
$ws.Range("M110").Value = 1655.5
$ws.Range("H132").Value = 1632
$ws.Range("I132").Value = 1548.75
$ws.Range("J132").Value = 1965
$ws.Range("K132").Value = 4646.25
$ws.Range("L132").Value = 5895
$ws.Range("M132").Value = -2116.25
$ws.Range("N132").Value = -10955

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 29700.857
$ws.Range("J82").Value = 41819.375
$ws.Range("L82").Value = 41819.375
$ws.Range("N82").Value = -42585.375
$ws.Range("H85").Value = 29700.857
$ws.Range("J85").Value = 41819.375
$ws.Range("L85").Value = 41819.375
$ws.Range("N85").Value = -44471.375
$ws.Range("H134").Value = 6352.8335
$ws.Range("I134").Value = 6450.905
$ws.Range("K134").Value = 19352.715
$ws.Range("M134").Value = -16817.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 52.88889
$ws.Range("I7").Value = 49.857143
$ws.Range("K7").Value = 49.857143
$ws.Range("M7").Value = 63.142857
$ws.Range("H17").Value = 3603.889
$ws.Range("I17").Value = 2310.25
$ws.Range("J17").Value = 13953
$ws.Range("K17").Value = 2310.25
$ws.Range("L17").Value = 13953
$ws.Range("M17").Value = -2136.25
$ws.Range("N17").Value = -14301
$ws.Range("H31").Value = 3139.8333
$ws.Range("I31").Value = 2753
$ws.Range("J31").Value = 3913.5
$ws.Range("K31").Value = 2753
$ws.Range("L31").Value = 3913.5
$ws.Range("M31").Value = -2458
$ws.Range("N31").Value = -4503.5
$ws.Range("H34").Value = 3139.8333
$ws.Range("I34").Value = 2753
$ws.Range("J34").Value = 3913.5
$ws.Range("K34").Value = 2753
$ws.Range("L34").Value = 3913.5
$ws.Range("M34").Value = -2551
$ws.Range("N34").Value = -4317.5
$ws.Range("H41").Value = 16829.3
$ws.Range("I41").Value = 2766.3333
$ws.Range("J41").Value = 22856.285
$ws.Range("K41").Value = 2766.3333
$ws.Range("L41").Value = 22856.285
$ws.Range("M41").Value = -2338.3333
$ws.Range("N41").Value = -23712.285
$ws.Range("H58").Value = 2946.5
$ws.Range("I58").Value = 2519.4167
$ws.Range("J58").Value = 3587.125
$ws.Range("K58").Value = 2519.4167
$ws.Range("L58").Value = 3587.125
$ws.Range("M58").Value = -2316.4167
$ws.Range("N58").Value = -3993.125
$ws.Range("H59").Value = 32754.889
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 34974.25
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 34974.25
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -37264.25
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H136").Value = 2946.5
$ws.Range("I136").Value = 2519.4167
$ws.Range("J136").Value = 3587.125
$ws.Range("K136").Value = 7558.250100000001
$ws.Range("L136").Value = 10761.375
$ws.Range("M136").Value = -5008.250100000001
$ws.Range("N136").Value = -15861.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 934.3333
$ws.Range("J32").Value = 934.3333
$ws.Range("L32").Value = 2802.9999
$ws.Range("N32").Value = -3368.9999
$ws.Range("H92").Value = 408.5
$ws.Range("I92").Value = 408.5
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1225.5
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 22.5
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H64").Value = 28633
$ws.Range("J64").Value = 28633
$ws.Range("L64").Value = 28633
$ws.Range("N64").Value = -29083
$ws.Range("H67").Value = 28633
$ws.Range("J67").Value = 28633
$ws.Range("L67").Value = 28633
$ws.Range("N67").Value = -30193
$ws.Range("H70").Value = 55000
$ws.Range("J70").Value = 55000
$ws.Range("L70").Value = 55000
$ws.Range("N70").Value = -55540
$ws.Range("H73").Value = 55000
$ws.Range("J73").Value = 55000
$ws.Range("L73").Value = 55000
$ws.Range("N73").Value = -56872
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3098.1304
$ws.Range("I122").Value = 2827.7693
$ws.Range("K122").Value = 8483.3079
$ws.Range("M122").Value = -6033.3079

Write-Output "done"